# Update the Sprites table: add a "Type" row, a "defaultZ" row, drop the
# platform1/platform2/Moveable rows/cols, rename player1 -> player, and add
# a new rock1/r1 entity column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any existing contents first so stale cells (e.g. old H column data)
# don't linger if the new range is a different shape.
$ws.Cells.Clear()

$headers = @("Name", "player", "enemy1", "enemy2", "block1", "block2", "staircase1", "staircase2", "rock1")
$ids     = @("ID",   "p1",     "e1",     "e2",     "b1",     "b2",     "s1",         "s2",         "r1")
$types   = @("Type", "entity", "entity", "entity", "block",  "block",  "block",      "block",      "block")

for ($col = 1; $col -le 9; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
    $ws.Cells.Item(2, $col).Value = $ids[$col - 1]
}

$ws.Cells.Item(3, 1).Value = "Solid"
for ($col = 2; $col -le 9; $col++) {
    $ws.Cells.Item(3, $col).Value = $true
}

for ($col = 1; $col -le 9; $col++) {
    $ws.Cells.Item(4, $col).Value = $types[$col - 1]
}

$ws.Cells.Item(5, 1).Value = "defaultZ"
for ($col = 2; $col -le 9; $col++) {
    $ws.Cells.Item(5, $col).Value = 2
}

# Restore the selection the author left the workbook in.
$ws.Range("A2").Select()
